$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking values stay stored as text (matches source sharedStrings)
$ws.Range("C2:D2").NumberFormat = "@"

# Update the data row (row 2) with the new values
$ws.Range("A2").Value = "dada life"
$ws.Range("B2").Value = "violents1"
$ws.Range("C2").Value = "2321"
$ws.Range("D2").Value = "88"
$ws.Range("E2").Value = "31/10/2022 20:42:56"

# Restore the default (General) cell style for C2:D2 now that the text values
# are locked in, so no extra style entries linger on these cells
$ws.Range("C2:D2").ClearFormats()

# Auto-fit columns A:E so the stored widths match the new content
$ws.Columns("A:E").AutoFit() | Out-Null
